$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data row was recorded. It was inserted at row 11 (the data
# that previously occupied rows 11-169 moves down to rows 12-170), carrying
# along its original formatting/values, and then the date (D) and volume (J)
# cells of the newly inserted row are set to the new observation.
$ws.Rows.Item(11).Copy()
$ws.Rows.Item(11).Insert()

$ws.Range("D11").Value = 44496
$ws.Range("J11").Value = 2900
